{"js": "// Replace the 24 division-problem text runs in the practice table with\n// their new values, per the commit's regenerated numbers.\nconst replacements = [\n  [\"457\u00f74=\", \"832\u00f75=\"],\n  [\"765\u00f74=\", \"187\u00f72=\"],\n  [\"530\u00f73=\", \"129\u00f72=\"],\n  [\"401\u00f72=\", \"743\u00f76=\"],\n  [\"450\u00f76=\", \"133\u00f74=\"],\n  [\"920\u00f79=\", \"808\u00f78=\"],\n  [\"573\u00f78=\", \"993\u00f78=\"],\n  [\"776\u00f75=\", \"330\u00f77=\"],\n  [\"555\u00f79=\", \"883\u00f73=\"],\n  [\"421\u00f73=\", \"701\u00f77=\"],\n  [\"772\u00f75=\", \"492\u00f74=\"],\n  [\"641\u00f73=\", \"499\u00f77=\"],\n  [\"934\u00f72=\", \"674\u00f75=\"],\n  [\"222\u00f74=\", \"756\u00f75=\"],\n  [\"781\u00f77=\", \"488\u00f78=\"],\n  [\"379\u00f78=\", \"588\u00f79=\"],\n  [\"534\u00f75=\", \"507\u00f73=\"],\n  [\"810\u00f78=\", \"726\u00f79=\"],\n  [\"938\u00f76=\", \"878\u00f72=\"],\n  [\"717\u00f75=\", \"395\u00f72=\"],\n  [\"579\u00f74=\", \"656\u00f76=\"],\n  [\"869\u00f74=\", \"832\u00f74=\"],\n  [\"823\u00f76=\", \"259\u00f79=\"],\n  [\"243\u00f76=\", \"410\u00f74=\"],\n  [\"402\u00f73=\", \"594\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 24 division-problem text runs in the practice table with\n# their new values, per the commit's regenerated numbers.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"457\u00f74=\", \"832\u00f75=\"),\n    @(\"765\u00f74=\", \"187\u00f72=\"),\n    @(\"530\u00f73=\", \"129\u00f72=\"),\n    @(\"401\u00f72=\", \"743\u00f76=\"),\n    @(\"450\u00f76=\", \"133\u00f74=\"),\n    @(\"920\u00f79=\", \"808\u00f78=\"),\n    @(\"573\u00f78=\", \"993\u00f78=\"),\n    @(\"776\u00f75=\", \"330\u00f77=\"),\n    @(\"555\u00f79=\", \"883\u00f73=\"),\n    @(\"421\u00f73=\", \"701\u00f77=\"),\n    @(\"772\u00f75=\", \"492\u00f74=\"),\n    @(\"641\u00f73=\", \"499\u00f77=\"),\n    @(\"934\u00f72=\", \"674\u00f75=\"),\n    @(\"222\u00f74=\", \"756\u00f75=\"),\n    @(\"781\u00f77=\", \"488\u00f78=\"),\n    @(\"379\u00f78=\", \"588\u00f79=\"),\n    @(\"534\u00f75=\", \"507\u00f73=\"),\n    @(\"810\u00f78=\", \"726\u00f79=\"),\n    @(\"938\u00f76=\", \"878\u00f72=\"),\n    @(\"717\u00f75=\", \"395\u00f72=\"),\n    @(\"579\u00f74=\", \"656\u00f76=\"),\n    @(\"869\u00f74=\", \"832\u00f74=\"),\n    @(\"823\u00f76=\", \"259\u00f79=\"),\n    @(\"243\u00f76=\", \"410\u00f74=\"),\n    @(\"402\u00f73=\", \"594\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n}\n"}
